$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "MEC-3B-Calderaria"

$ws.Range("D3").Value = "MEC-3B-Calderaria"
$ws.Range("E3").Value = "-"

$ws.Range("D4").Value = "MEC-3B-Calderaria"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "MCT-1A-Metrologia"

$ws.Range("D6").Value = "MEC-3B-Calderaria"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "MCT-1A-Metrologia"

$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "MCT-1A-Metrologia"

$ws.Range("F8").Value = "MCT-1A-Metrologia"

$ws.Range("B10").Value = "MEC-3A-Calderaria"
$ws.Range("E10").Value = "MEC-3A-Calderaria"

$ws.Range("E11").Value = "MEC-3A-Calderaria"

$ws.Range("D12").Value = "-"

$ws.Range("D14").Value = "-"

$ws.Range("D15").Value = "-"

$ws.Range("C16").Value = "MEC-3A-Calderaria"
$ws.Range("D16").Value = "-"
